# Update odds values in Sheet1 according to the upstream FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.8
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1.17
$ws.Range("K2").Value = 5
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 13
$ws.Range("X2").Value = 34
$ws.Range("AG2").Value = 29

# Row 3
$ws.Range("G3").Value = 1.57
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 6.4
$ws.Range("K3").Value = 5.8
$ws.Range("R3").Value = 2.42
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 4.7
$ws.Range("U3").Value = 5.8
$ws.Range("W3").Value = 10.5
$ws.Range("Z3").Value = 5.8
$ws.Range("AA3").Value = 7.2
$ws.Range("AC3").Value = 200
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 37
$ws.Range("AH3").Value = 100

# Row 8
$ws.Range("N8").Value = 1.93
$ws.Range("O8").Value = 1.88

# Row 9
$ws.Range("G9").Value = 1.52
$ws.Range("H9").Value = 4.4
$ws.Range("I9").Value = 4.85
$ws.Range("N9").Value = 1.45
$ws.Range("O9").Value = 2.37
$ws.Range("T9").Value = 10
$ws.Range("AD9").Value = 19

# Row 11
$ws.Range("W11").Value = 11.5

# Row 12
$ws.Range("N12").Value = 2.15
$ws.Range("O12").Value = 1.62
$ws.Range("AB12").Value = 14.5

# Row 16
$ws.Range("AD16").Value = 9.5
$ws.Range("AE16").Value = 19
